# Checklist.xlsx update - mark additional checklist items as completed / reset a few others,
# and move the active selection down to E25.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Mark these checklist rows as completed ("Y")
$ws.Range("E5").Value = "Y"
$ws.Range("E6").Value = "Y"
$ws.Range("E7").Value = "Y"
$ws.Range("E8").Value = "Y"
$ws.Range("E9").Value = "Y"
$ws.Range("E10").Value = "Y"
$ws.Range("E11").Value = "Y"
$ws.Range("E13").Value = "Y"
$ws.Range("E14").Value = "Y"
$ws.Range("E24").Value = "Y"

# Clear these checklist marks (no longer completed / remove the tentative "Y?")
$ws.Range("E19").ClearContents()
$ws.Range("E20").ClearContents()
$ws.Range("E21").ClearContents()
$ws.Range("E22").ClearContents()

# Move the active selection/cursor to E25
$ws.Range("E25").Select()
